$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.278.99"
$ws.Range("E2").Value = "  +5.11%  "
$ws.Range("D3").Value = "1.911.66"
$ws.Range("E3").Value = "  +5.40%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9999"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "252.76"
$ws.Range("E5").Value = "  +0.74%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9996"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5140"
$ws.Range("E7").Value = "  +3.14%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "45.48"
$ws.Range("E8").Value = "  +5.16%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2979"
$ws.Range("E9").Value = "  +7.15%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06787"
$ws.Range("E10").Value = "  +6.31%  "
$ws.Range("D11").Value = "1.914.69"
$ws.Range("E11").Value = "  +5.57%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "17.44"
$ws.Range("E12").Value = "  +3.86%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07402"
$ws.Range("E13").Value = "  +3.36%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6959"
$ws.Range("E14").Value = "  +7.10%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "87.47"
$ws.Range("E15").Value = "  +6.60%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.879"
$ws.Range("E16").Value = "  +3.58%  "
$ws.Range("D17").Value = "30.277.99"
$ws.Range("E17").Value = "  +5.19%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008075"
$ws.Range("E18").Value = "  +9.25%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9999"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.00"
$ws.Range("E20").Value = "  +5.83%  "
$ws.Range("D21").Value = "2.162.14"
$ws.Range("E21").Value = "  +5.62%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9995"
$ws.Range("E22").Value = "  +0.10%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.844"
$ws.Range("E23").Value = "  +4.99%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.752"
$ws.Range("E24").Value = "  +7.44%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.155"
$ws.Range("E25").Value = "  +2.93%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "146.27"
$ws.Range("E26").Value = "  +1.58%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "137.72"
$ws.Range("E27").Value = "  +15.66%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.21"
$ws.Range("E28").Value = "  +7.21%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.012"
$ws.Range("E29").Value = "  +6.17%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.400"
$ws.Range("E30").Value = "  +0.14%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.259"
$ws.Range("E31").Value = "  +1.71%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08823"
$ws.Range("E32").Value = "  +5.44%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.009"
$ws.Range("E33").Value = "  +3.96%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05115"
$ws.Range("E34").Value = "  +2.98%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.156"
$ws.Range("E35").Value = "  +6.07%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7202"
$ws.Range("E36").Value = "  +5.76%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.689"
$ws.Range("E37").Value = "  +0.23%  "
$ws.Range("B38").Value = "MXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.835"
$ws.Range("E38").Value = "  +3.35%  "
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.307"
$ws.Range("E39").Value = "  +5.14%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9745"
$ws.Range("E40").Value = "  +0.55%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01697"
$ws.Range("E41").Value = "  +6.79%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.093"
$ws.Range("E42").Value = "  +1.39%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4316"
$ws.Range("E43").Value = "  +4.78%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "105.74"
$ws.Range("E44").Value = "  +4.15%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9989"
$ws.Range("E45").Value = "  -0.07%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.689"
$ws.Range("E46").Value = "  +6.21%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1282"
$ws.Range("E47").Value = "  +4.59%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05756"
$ws.Range("E48").Value = "  +4.60%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "33.27"
$ws.Range("E49").Value = "  +4.91%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.500"
$ws.Range("E50").Value = "  +3.91%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3824"
$ws.Range("E51").Value = "  +4.66%  "
